# Update sheet data rows 2-31 (replace old scraped listings with new
# Cloudscraper-sourced Zonaprop results) and extend the used range.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('Balvanera, Capital Federal', '$ 700.000', '41 m² tot.3 amb.2 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-alquiler-3-ambientes-dueno-caba-balvanera-58080275.html?n_src=Listado&n_pills=Encargado&n_pg=1&n_pos=1')
    ,@('Almagro Sur, Almagro', '$ 800.000', '42 m² tot.2 amb.1 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquiler-2-amb-42-m-sup2--excelente-estado.-muy-58093867.html?n_src=Listado&n_pills=Lavadero&n_pg=1&n_pos=2')
    ,@('Monserrat, Capital Federal', '$ 1.400.000', '41 m² tot.3 amb.2 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-de-3-ambientes-amoblado-listo-para-58096171.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=3')
    ,@('Tribunales, Capital Federal', '$ 680.000', '37 m² tot.2 amb.1 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-depto-equipado-y-funcional-en-inmejorable-ubicacion!-52711772.html?n_src=Listado&n_pills=Permite+mascotas&n_pg=1&n_pos=4')
    ,@('Palermo Soho, Palermo', '$ 600.000', '65 m² tot.1 amb.1 dorm.1 baño1 coch.', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquilo-departamento-totalmente-amueblado-58096033.html?n_src=Listado&n_pills=Permite+mascotas&n_pg=1&n_pos=5')
    ,@('Caballito Sur, Caballito', '$ 1.200.000', '83 m² tot.4 amb.3 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-excelente-dpto-4-amb-c-dep-al-frente-58096013.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=6')
    ,@('Recoleta, Capital Federal', '$ 450.000', '29 m² tot.2 amb.1 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-dueno-alquila-dpto-2-ambientes-lateral-super-luminoso-57250660.html?n_src=Listado&n_pills=Lavadero&n_pg=1&n_pos=7')
    ,@('Almagro Norte, Almagro', '$ 650.000', '38 m² tot.2 amb.1 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-dueno-alquila-dpto-2-amb-vista-abierta-super-luminoso-57652927.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=8')
    ,@('Villa del Parque, Capital Federal', '$ 750.000', '41 m² tot.2 amb.1 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-con-terraza-y-sol-57938037.html?n_src=Listado&n_pills=Encargado&n_pg=1&n_pos=9')
    ,@('Belgrano R, Belgrano', '$ 700.000', '60 m² tot.2 amb.1 dorm.1 baño1 coch.', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-en-belgrano-r-dueno-directo-58095968.html?n_src=Listado&n_pills=Pileta&n_pg=1&n_pos=10')
    ,@('Flores Norte, Flores', '$ 620.000', '75 m² tot.3 amb.2 dorm.1 baño1 coch.', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-3-ambientes-58095956.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=11')
    ,@('Flores Norte, Flores', '$ 620.000', '75 m² tot.3 amb.2 dorm.1 baño1 coch.', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-3-ambientes-58095955.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=12')
    ,@('Floresta Sur, Floresta', '$ 750', '89 m² tot.3 amb.2 dorm.2 baños1 coch.', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquilo-depto-amueblado-58095949.html?n_src=Listado&n_pills=Permite+mascotas&n_pg=1&n_pos=13')
    ,@('Villa Crespo, Capital Federal', 'USD 800', '72 m² tot.3 amb.2 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquiler-dpto-3-ambientes-dueno-directo-amoblado-58095947.html?n_src=Listado&n_pills=Permite+mascotas&n_pg=1&n_pos=14')
    ,@('Flores Norte, Flores', '$ 620.000', '75 m² tot.3 amb.2 dorm.1 baño1 coch.', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-3-ambientes-58095948.html?n_src=Listado&n_pills=Permite+mascotas&n_pg=1&n_pos=15')
    ,@('Flores Norte, Flores', '$ 620.000', '75 m² tot.3 amb.2 dorm.1 baño1 coch.', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-3-ambientes-58095946.html?n_src=Listado&n_pills=Laundry&n_pg=1&n_pos=16')
    ,@('Flores Norte, Flores', '$ 630.000', '70 m² tot.2 amb.1 dorm.1 baño1 coch.', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-2-ambiente-58095945.html?n_src=Listado&n_pills=Encargado&n_pg=1&n_pos=17')
    ,@('Flores Norte, Flores', '$ 630.000', '75 m² tot.3 amb.2 dorm.1 baño1 coch.', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-2-ambiente-58095921.html?n_src=Listado&n_pills=Permite+mascotas&n_pg=1&n_pos=18')
    ,@('Flores Norte, Flores', '$ 620.000', '75 m² tot.3 amb.2 dorm.1 baño1 coch.', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-3-ambientes-58095925.html?n_src=Listado&n_pills=Laundry&n_pg=1&n_pos=19')
    ,@('Villa Urquiza, Capital Federal', '$ 850.000', '48 m² tot.1 amb.1 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-2-ambientes-en-alquiler-al-frente-villa-urquiza-dueno-58095610.html?n_src=Listado&n_pills=Laundry&n_pg=1&n_pos=20')
    ,@('San Cristobal, Capital Federal', '$ 499.000', '37 m² tot.1 amb.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquiler-monoambiente.-bajas-expensas.-sin-comision.-54758980.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=21')
    ,@('Villa General Mitre, Capital Federal', '$ 700.000', '39 m² tot.2 amb.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquiler-dos-ambientes-villa-general-mitre-58095085.html?n_src=Listado&n_pills=Pileta&n_pg=1&n_pos=22')
    ,@('Parque Rivadavia, Caballito', '$ 600.000', '42 m² tot.2 amb.1 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-dos-ambientes-en-caballito.-58094336.html?n_src=Listado&n_pills=Lavadero&n_pg=1&n_pos=23')
    ,@('San Cristobal, Capital Federal', '$ 650.000', '35 m² tot.', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquilo-departamento-dos-ambientes-san-cristobal-58094313.html?n_src=Listado&n_pg=1&n_pos=24')
    ,@('Recoleta, Capital Federal', 'USD 700', '200 m² tot.4 amb.3 dorm.2 baños', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-dueno-directo-todo-incluido-resuelvo-hoy.-video-58094288.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=25')
    ,@('San Cristobal, Capital Federal', '$ 980.000', '74 m² tot.3 amb.2 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-departamento-3-ambientes-luminoso-dueno-directo-58094293.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=26')
    ,@('Belgrano, Capital Federal', '$ 950.000', '52 m² tot.2 amb.1 dorm.1 baño1 coch.', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-duena-alquila-2-ambientes-con-cochera.-amplio-y-52962262.html?n_src=Listado&n_pills=Parrilla&n_pg=1&n_pos=27')
    ,@('Lomas de Núñez, Núñez', 'USD 750', '53 m² tot.2 amb.1 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquiler-departamento-de-2-ambientes-en-nunez-caba-58094208.html?n_src=Listado&n_pills=Parrilla&n_pg=1&n_pos=28')
    ,@('Barracas, Capital Federal', '$ 700.000', '65 m² tot.3 amb.2 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-dueno-alquila-luminoso-3-amb.-en-barracas-9-piso-al-57960079.html?n_src=Listado&n_pills=Lavadero&n_pg=1&n_pos=29')
    ,@('Almagro, Capital Federal', '$ 480.000', '32 m² tot.2 amb.1 dorm.1 baño', 'https://www.zonaprop.com.ar/propiedades/clasificado/alclappa-alquiler-departamento-dos-ambientes.-mario-bravo-al-58093565.html?n_src=Listado&n_pills=Aire+acondicionado&n_pg=1&n_pos=30')
)

# Column B mixes "$ ###.###" / "USD ###" price strings. Values that look
# numeric (leading "$") get auto-parsed as currency by Excel unless the
# cell is pre-formatted as Text, so force Text format on those rows first
# to keep the literal price string (matches the source report).
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $price = $data[$i][1]
    if ($price.StartsWith("$")) {
        $ws.Cells.Item($row, 2).NumberFormat = "@"
    }
}

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

